$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected cells to be treated as plain text so that
# numeric-looking values (e.g. "138.15") are not reinterpreted as numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '63.143.95'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").Value = '3.138.29'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '589.94'
$ws.Range("E5").Value = '  -1.66%  '
$ws.Range("D6").Value = '138.15'
$ws.Range("E6").Value = '  -3.06%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.136.43'
$ws.Range("E8").Value = '  +0.49%  '
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("E10").Value = '  -1.65%  '
$ws.Range("D11").Value = '5.29'
$ws.Range("E11").Value = '  -1.26%  '
$ws.Range("E12").Value = '  -1.87%  '
$ws.Range("E13").Value = '  -2.76%  '
$ws.Range("D14").Value = '34.18'
$ws.Range("E14").Value = '  -2.41%  '
$ws.Range("D15").Value = '3.663.20'
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("D17").Value = '3.141.30'
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("D18").Value = '63.089.93'
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("E19").Value = '  -2.35%  '
$ws.Range("D20").Value = '473.90'
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("D21").Value = '14.07'
$ws.Range("E21").Value = '  -3.06%  '
$ws.Range("D22").Value = '0.698'
$ws.Range("E22").Value = '  -1.05%  '
$ws.Range("E23").Value = '  +0.79%  '
$ws.Range("D24").Value = '84.57'
$ws.Range("E24").Value = '  -3.61%  '
$ws.Range("D25").Value = '13.02'
$ws.Range("E25").Value = '  -2.30%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").Value = '2.71'
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("E28").Value = '  -3.40%  '
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").Value = '2.11'
$ws.Range("E29").Value = '  +3.06%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").Value = '6.94'
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("D32").Value = '26.81'
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("E33").Value = '  -5.42%  '
$ws.Range("D34").Value = '2.53'
$ws.Range("E34").Value = '  -4.31%  '
$ws.Range("E35").Value = '  -1.60%  '
$ws.Range("D36").Value = '5.80'
$ws.Range("E36").Value = '  -3.10%  '
$ws.Range("D37").Value = '52.36'
$ws.Range("E37").Value = '  -0.45%  '
$ws.Range("D38").Value = '0.0₃0696'
$ws.Range("E38").Value = '  -7.04%  '
$ws.Range("D39").Value = '0.0387'
$ws.Range("E39").Value = '  -0.84%  '
$ws.Range("D40").Value = '421.12'
$ws.Range("E40").Value = '  -3.36%  '
$ws.Range("E41").Value = '  -6.87%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D43").Value = '2.909.77'
$ws.Range("E43").Value = '  +1.52%  '
$ws.Range("E44").Value = '  -5.37%  '
$ws.Range("D45").Value = '0.262'
$ws.Range("E45").Value = '  +1.76%  '
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").Value = '2.13'
$ws.Range("E47").Value = '  -3.82%  '
$ws.Range("D48").Value = '25.39'
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  -7.06%  '
$ws.Range("D51").Value = '120.47'
$ws.Range("E51").Value = '  -0.64%  '

# Restore the original (default) cell style now that values are set.
$ws.Range("B2:E51").Style = "Normal"
